$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $r.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $r.Text = $new
}

function Replace-Range-Between($startText, $endText, $new) {
    # Finds startText, then finds endText after it, and replaces the whole
    # combined span (startText ... endText) with $new.
    $r1 = $d.Content
    $r1.Find.Execute($startText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $r2 = $d.Content
    $r2.Start = $r1.End
    $r2.Find.Execute($endText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $r1.End = $r2.End
    $r1.Text = $new
}

# --- Title / byline / email ---
Replace-Text "Unraveling the Enigma of Dark Matter" "The Art of Self-Expression through Creative Writing"
Replace-Text "Alexander Forsyth" "Jane Williams"
Replace-Text "alexander" "jane"
Replace-Text "forsyth@cosmology" "williams@academyhigh"
Replace-Text "edu" "org"

# --- Body paragraph 1 ---
Replace-Text "Like a cosmic puzzle, dark matter has captivated the imagination of scientists for decades" "The universe of literature is an avenue of uncharted spaces where untold tales await discovery, similar to the first stargazers sending telescopes into the cosmos"
Replace-Text " This elusive substance, believed to constitute over 80% of the universe's mass, remains shrouded in mystery" " It is through the lens of words that we decipher our innermost thoughts, much like a botanist identifies a species through its anatomy"

Replace-Range-Between " Its existence is inferred through its gravitational influence on visible matter, yet its true nature and properties remain enigmatic" " As we unravel the secrets of dark matter, we uncover new insights into the fundamental fabric of our universe" " The careful articulation of language mirrors biological taxonomy, as each word carefully dissects and categorizes the fragments of our thoughts, evolving them into cohesive ideas"

Replace-Text "A multitude of astrophysical observations provide compelling evidence for the existence of dark matter" "The act of creative writing is like traversing a maze, where the first step is selecting the appropriate narrative form--be it poetry's focused energy or the far-reaching epic of a novel"

Replace-Range-Between " Studies of galaxy rotation curves reveal an unexpected distribution of mass, with the outer regions of galaxies rotating faster than predicted by the visible mass alone" " Moreover, the cosmic microwave background radiation, a relic of the early universe, bears the imprint of dark matter's influence, providing a glimpse into its properties and distribution" " Once the narrative space is secured, one must navigate the sprawling landscape of language, where each word is meticulously chosen, akin to a master chef harmonizing ingredients to create a symphony of flavors"

Replace-Text "Theoretical models offer frameworks for understanding the nature of dark matter" "Language itself is a prism through which human emotions are deciphered, much like a scholar decodes ancient inscriptions"

Replace-Range-Between " One prominent hypothesis is that dark matter consists of weakly interacting massive particles (WIMPs), subatomic particles with masses ranging from 10 to 1,000 times that of the proton" " Alternatively, modified gravity theories propose that the observed effects attributed to dark matter may arise from modifications to the laws of gravity on large scales" " Just as archeologists exhume fragments of an ancient civilization, writers unearth their memories, hopes, and dreams to weave them into intricate tapestries of language, allowing readers to appreciate the exquisite beauty of their humanity"

# --- Summary paragraph ---
Replace-Text "Dark matter, a mysterious and elusive substance, continues to challenge our understanding of the universe" "Creative writing is a unique blend of self-expression and artistic endeavor, mirroring various academic disciplines"
Replace-Text " Through meticulous observations and theoretical exploration, scientists are unraveling the enigma of dark matter, piece by piece" " Writers explore the terrain of their inner selves, akin to explorers venturing into uncharted territories, merging thoughtful selection of narrative form with meticulous word choice"

Replace-Range-Between " The quest " "to unveil its true nature is not merely an academic pursuit; it holds the key to unlocking fundamental questions about the cosmos, from its origin and evolution to its ultimate fate" " Through this process, writers shed new light on essential human emotions by excavating memories from subconscious archives, revealing the depths of human existence"

Replace-Text " As we delve deeper into the mysteries of dark matter, we expand the boundaries of human knowledge and gain a profound appreciation for the intricate tapestry of the universe we inhabit" " Creative writing is a potent art form that allows people to comprehend and express their innermost realities"

# --- Add a new empty paragraph at the very end of the document body ---
$endPos = $d.Content.End
$endRange = $d.Range($endPos, $endPos)
$endRange.Text = "`r"
